# Generate Report for Handoff
# - Row 2 (the 79efc3a7... file) moves from "Handed back: in sync with en-US"
#   to "Ready for handoff", with refreshed handoff timestamps.
# - Row 3 (the ea42df71... file) is dropped entirely from every sheet, along
#   with its hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-38-19 02:38:21"

# Drop row 3 (ea42df71-...) entirely - shifts dimension from A1:D3 to A1:D2
$ws1.Rows(3).Delete()

# Rebuild hyperlinks so the orphaned ea42df71 link (old A3) is gone and the
# surviving one keeps pointing at the same target as before.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0fc54bcd55584786bde15e589d5eaccb790d253/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md")
$ws1.Range("A2").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-19 02:38:19"

# Drop row 3 (ea42df71-...) entirely - shifts dimension from A1:K3 to A1:K2
$ws2.Rows(3).Delete()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0fc54bcd55584786bde15e589d5eaccb790d253/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0fc54bcd55584786bde15e589d5eaccb790d253/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/200a34eaf1362ed1a043d104f105b48c078d17e3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5b400ac3a4a197a0071e5196c625eae27716988f/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa3a0eed628e0b24b73d74959a6e7ccb9c9c1ee3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.zh-cn.xlf")
$ws2.Range("A2,B2,D2,F2,G2").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-19 02:38:21"

# Drop row 3 (ea42df71-...) entirely - shifts dimension from A1:K3 to A1:K2
$ws3.Rows(3).Delete()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0fc54bcd55584786bde15e589d5eaccb790d253/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0fc54bcd55584786bde15e589d5eaccb790d253/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ffebd461be235cdf52406fa73d9fdbdcd0d13c73/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b804122ea0b457b2932cfab425cc6ee5ca18e3d5/e2e/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/47eb3f961a8c8a74a9c64a8c1bffd7edfcbb9305/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf", "", "", "79efc3a7-98d0-4f6b-a706-c7388cd3ac8a.6a0d9d2bb57d903e710272a8ba7ce0f4e005fb66.de-de.xlf")
$ws3.Range("A2,B2,D2,F2,G2").Style = "HyperLink"

Write-Host "Handoff report regenerated: dropped ea42df71 row, refreshed 79efc3a7 status/timestamps."
